# Append 34 new store rows (182-215) to Sheet1, matching the source workbook upload.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(182, 1).Value = 'Do or Dice '
$ws.Cells.Item(182, 2).Value = 'Unit 6b, One, Market St, Addlestone KT15 2GB'
$ws.Cells.Item(182, 3).Value = 'KT15 2GB'
$ws.Cells.Item(182, 5).Value = -0.48982999999999999
$ws.Cells.Item(182, 6).Value = 51.370089999999998

$ws.Cells.Item(183, 1).Value = 'The Missing Geek '
$ws.Cells.Item(183, 2).Value = '11 Clarence St, Staines TW18 4SU'
$ws.Cells.Item(183, 3).Value = 'TW18 4SU'
$ws.Cells.Item(183, 4).Value = 1784557007
$ws.Cells.Item(183, 5).Value = -0.51437940000000004
$ws.Cells.Item(183, 6).Value = 51.433722299999999

$ws.Cells.Item(184, 1).Value = 'Collectors Cardhouse'
$ws.Cells.Item(184, 2).Value = 'Lower Rd, Chorleywood, Rickmansworth WD3 5LH'
$ws.Cells.Item(184, 3).Value = 'WD3 5LH'
$ws.Cells.Item(184, 5).Value = -0.52002769999999998
$ws.Cells.Item(184, 6).Value = 51.654079099999997

$ws.Cells.Item(185, 1).Value = 'Tabletop Republic - High Wycombe'
$ws.Cells.Item(185, 2).Value = '37 High St, High Wycombe HP11 2AG'
$ws.Cells.Item(185, 3).Value = 'HP11 2AG'
$ws.Cells.Item(185, 4).Value = 1494257566
$ws.Cells.Item(185, 5).Value = -0.75032410000000005
$ws.Cells.Item(185, 6).Value = 51.628897500000001

$ws.Cells.Item(186, 1).Value = 'Axion Now'
$ws.Cells.Item(186, 2).Value = 'Merritt House, 1 Hill Ave, Amersham HP6 5BQ'
$ws.Cells.Item(186, 3).Value = 'HP6 5BQ'
$ws.Cells.Item(186, 4).Value = 1494728181
$ws.Cells.Item(186, 5).Value = -0.60797389999999996
$ws.Cells.Item(186, 6).Value = 51.674631499999997

$ws.Cells.Item(187, 1).Value = 'Eclectic Games'
$ws.Cells.Item(187, 2).Value = '5 Union St, Reading RG1 1EU'
$ws.Cells.Item(187, 3).Value = 'RG1 1EU'
$ws.Cells.Item(187, 4).Value = 1189598250
$ws.Cells.Item(187, 5).Value = -0.9734507
$ws.Cells.Item(187, 6).Value = 51.456302100000002

$ws.Cells.Item(188, 1).Value = 'Elemental Cards'
$ws.Cells.Item(188, 2).Value = 'Holme Grange Craft Village, Unit 4 Heathlands Rd, Wokingham RG40 3AW'
$ws.Cells.Item(188, 3).Value = 'RG40 3AW'
$ws.Cells.Item(188, 4).Value = 7908441130
$ws.Cells.Item(188, 5).Value = -0.81596999999999997
$ws.Cells.Item(188, 6).Value = 51.39414

$ws.Cells.Item(189, 1).Value = 'Osv Gaming'
$ws.Cells.Item(189, 2).Value = '25a waterbourne walk, Leighton Buzzard LU7 1DH'
$ws.Cells.Item(189, 4).Value = '01525 377716'
$ws.Cells.Item(189, 3).Value = 'LU7 1DH'
$ws.Cells.Item(189, 5).Value = -0.66358289999999998
$ws.Cells.Item(189, 6).Value = 51.916846200000002

$ws.Cells.Item(190, 1).Value = 'Clockwork Dragon'
$ws.Cells.Item(190, 2).Value = '42 Bath St, Abingdon OX14 3QH'
$ws.Cells.Item(190, 3).Value = 'OX14 3QH'
$ws.Cells.Item(190, 4).Value = 1235528462
$ws.Cells.Item(190, 5).Value = -1.2836194000000001
$ws.Cells.Item(190, 6).Value = 51.670870899999997

$ws.Cells.Item(191, 1).Value = 'Atomic Cards'
$ws.Cells.Item(191, 2).Value = 'Unit T18, The I O Centre, Hobley Dr, Swindon SN3 4JG'
$ws.Cells.Item(191, 3).Value = 'SN3 4JG'
$ws.Cells.Item(191, 4).Value = 1793376112
$ws.Cells.Item(191, 5).Value = -1.7520024999999999
$ws.Cells.Item(191, 6).Value = 51.5827253

$ws.Cells.Item(192, 1).Value = 'The Card Market'
$ws.Cells.Item(192, 2).Value = '5 Market Hall, THE ARCADE, Bedford MK40 1NS'
$ws.Cells.Item(192, 3).Value = 'MK40 1NS'
$ws.Cells.Item(192, 5).Value = -0.46767910000000001
$ws.Cells.Item(192, 6).Value = 52.1374371

$ws.Cells.Item(193, 1).Value = 'TabletopMK'
$ws.Cells.Item(193, 2).Value = '59 Stratford Rd, Wolverton, Milton Keynes MK12 5LT'
$ws.Cells.Item(193, 3).Value = 'MK12 5LT'
$ws.Cells.Item(193, 5).Value = -0.8148725
$ws.Cells.Item(193, 6).Value = 52.062492399999996

$ws.Cells.Item(194, 1).Value = 'Wargames Workshop MK '
$ws.Cells.Item(194, 2).Value = '18-19/Kingston Quarter/Winchester Circle, Milton Keynes MK10 0BA'
$ws.Cells.Item(194, 3).Value = 'MK10 0BA'
$ws.Cells.Item(194, 5).Value = -0.68980859999999999
$ws.Cells.Item(194, 6).Value = 52.035299299999998

$ws.Cells.Item(195, 1).Value = 'The Pokemon Plug'
$ws.Cells.Item(195, 2).Value = '1 St Mary''s St, Huntingdon PE29 3PE'
$ws.Cells.Item(195, 3).Value = 'PE29 3PE'
$ws.Cells.Item(195, 4).Value = 1480431142
$ws.Cells.Item(195, 5).Value = -0.1819954
$ws.Cells.Item(195, 6).Value = 52.328454100000002

$ws.Cells.Item(196, 1).Value = 'Picklestix TCG'
$ws.Cells.Item(196, 2).Value = '4 Fishers Yard, St. Neots PE19 2AG'
$ws.Cells.Item(196, 3).Value = 'PE19 2AG'
$ws.Cells.Item(196, 4).Value = 7710137981
$ws.Cells.Item(196, 5).Value = -0.27056409999999997
$ws.Cells.Item(196, 6).Value = 52.2268179

$ws.Cells.Item(197, 1).Value = 'Kingdom Gaming'
$ws.Cells.Item(197, 2).Value = 'The Ridings, Northampton NN1 2BA'
$ws.Cells.Item(197, 3).Value = 'NN1 2BA'
$ws.Cells.Item(197, 4).Value = 1604458503
$ws.Cells.Item(197, 5).Value = -0.89150700000000005
$ws.Cells.Item(197, 6).Value = 52.238549599999999

$ws.Cells.Item(198, 1).Value = 'Unified Cards'
$ws.Cells.Item(198, 2).Value = '35a Corporation St, Willow Pl, Corby NN17 1NQ'
$ws.Cells.Item(198, 3).Value = 'NN17 1NQ'
$ws.Cells.Item(198, 4).Value = 7504722951
$ws.Cells.Item(198, 5).Value = -0.70262190000000002
$ws.Cells.Item(198, 6).Value = 52.487569200000003

$ws.Cells.Item(199, 1).Value = 'WildPhire Collectables'
$ws.Cells.Item(199, 2).Value = '91 High St, Rushden NN10 0NZ'
$ws.Cells.Item(199, 3).Value = 'NN10 0NZ'
$ws.Cells.Item(199, 4).Value = 7445045102
$ws.Cells.Item(199, 5).Value = -0.59847260000000002
$ws.Cells.Item(199, 6).Value = 52.291762599999998

$ws.Cells.Item(200, 1).Value = 'Athena Games'
$ws.Cells.Item(200, 2).Value = '9, St Gregory''s Alley, Norwich NR2 1ER'
$ws.Cells.Item(200, 3).Value = 'NR2 1ER'
$ws.Cells.Item(200, 4).Value = 1603460910
$ws.Cells.Item(200, 5).Value = 1.2911577999999999
$ws.Cells.Item(200, 6).Value = 52.6304309

$ws.Cells.Item(201, 1).Value = 'The Rift'
$ws.Cells.Item(201, 2).Value = 'Rivergate Arcade, Viersen Platz, Peterborough PE1 1EL'
$ws.Cells.Item(201, 3).Value = 'PE1 1EL'
$ws.Cells.Item(201, 4).Value = 1733341007
$ws.Cells.Item(201, 5).Value = -0.242892
$ws.Cells.Item(201, 6).Value = 52.569692099999997

$ws.Cells.Item(202, 1).Value = 'GG Trading Cards and Collectables'
$ws.Cells.Item(202, 2).Value = 'Pod 12/13, Stonham Barns, Pettaugh road, Stonham Aspal, Stowmarket IP14 6AT'
$ws.Cells.Item(202, 3).Value = 'IP14 6AT'
$ws.Cells.Item(202, 4).Value = 7861925738
$ws.Cells.Item(202, 5).Value = 1.1374873000000001
$ws.Cells.Item(202, 6).Value = 52.189046500000003

$ws.Cells.Item(203, 1).Value = 'Rocket''s Hideout!'
$ws.Cells.Item(203, 2).Value = '64a Gowthorpe, Selby YO8 4ET'
$ws.Cells.Item(203, 3).Value = 'YO8 4ET'
$ws.Cells.Item(203, 4).Value = 1757713998
$ws.Cells.Item(203, 5).Value = -1.0719958000000001
$ws.Cells.Item(203, 6).Value = 53.7837429

$ws.Cells.Item(204, 1).Value = 'Planet JJ''s - Geekery'
$ws.Cells.Item(204, 2).Value = '11, Quarry Hill Parade, Tonbridge TN9 2HR'
$ws.Cells.Item(204, 3).Value = 'TN9 2HR'
$ws.Cells.Item(204, 5).Value = 0.27000059999999998
$ws.Cells.Item(204, 6).Value = 51.190068500000002

$ws.Cells.Item(205, 1).Value = 'Chaos Cards Tabletop Gaming Centre'
$ws.Cells.Item(205, 2).Value = '100 Sandgate Rd, Folkestone CT20 2BW'
$ws.Cells.Item(205, 3).Value = 'CT20 2BW'
$ws.Cells.Item(205, 4).Value = 1303255522
$ws.Cells.Item(205, 5).Value = 1.1755243
$ws.Cells.Item(205, 6).Value = 51.077955899999999

$ws.Cells.Item(206, 1).Value = 'Eclipse Gaming'
$ws.Cells.Item(206, 2).Value = '39 Railway St, Chatham ME4 4RH'
$ws.Cells.Item(206, 3).Value = 'ME4 4RH'
$ws.Cells.Item(206, 4).Value = 1634553555
$ws.Cells.Item(206, 5).Value = 0.52291889999999996
$ws.Cells.Item(206, 6).Value = 51.381970799999998

$ws.Cells.Item(207, 1).Value = 'Protect and Collect'
$ws.Cells.Item(207, 2).Value = '61 Pelham Rd S, Gravesend DA11 8QS'
$ws.Cells.Item(207, 3).Value = 'DA11 8QS'
$ws.Cells.Item(207, 4).Value = 7763981878
$ws.Cells.Item(207, 5).Value = 0.35454989999999997
$ws.Cells.Item(207, 6).Value = 51.434489900000003

$ws.Cells.Item(208, 1).Value = 'Cataclysm Games UK'
$ws.Cells.Item(208, 2).Value = '45 Canterbury St, Gillingham ME7 5TR'
$ws.Cells.Item(208, 3).Value = 'ME7 5TR'
$ws.Cells.Item(208, 4).Value = 1634852303
$ws.Cells.Item(208, 5).Value = 0.54405590000000004
$ws.Cells.Item(208, 6).Value = 51.385871700000003

$ws.Cells.Item(209, 1).Value = 'LT Gaming'
$ws.Cells.Item(209, 2).Value = 'LT Gaming, 41 Robertson St, Hastings TN34 1HL'
$ws.Cells.Item(209, 3).Value = 'TN34 1HL'
$ws.Cells.Item(209, 5).Value = 0.57810980000000001
$ws.Cells.Item(209, 6).Value = 50.855362999999997

$ws.Cells.Item(210, 1).Value = 'Level Up Games'
$ws.Cells.Item(210, 2).Value = '40 Palace St, Canterbury CT1 2DZ'
$ws.Cells.Item(210, 3).Value = 'CT1 2DZ'
$ws.Cells.Item(210, 4).Value = 1227785002
$ws.Cells.Item(210, 5).Value = 1.081637
$ws.Cells.Item(210, 6).Value = 51.281276499999997

$ws.Cells.Item(211, 1).Value = 'Famous Collectables'
$ws.Cells.Item(211, 2).Value = '19 Western Rd, Bexhill-on-Sea TN40 1DU'
$ws.Cells.Item(211, 3).Value = 'TN40 1DU'
$ws.Cells.Item(211, 4).Value = 1424215577
$ws.Cells.Item(211, 5).Value = 0.47249089999999999
$ws.Cells.Item(211, 6).Value = 50.840087799999999

$ws.Cells.Item(212, 1).Value = 'Galleon Games'
$ws.Cells.Item(212, 2).Value = '30 St Leonards Rd, Bexhill-on-Sea TN40 1HT'
$ws.Cells.Item(212, 3).Value = 'TN40 1HT'
$ws.Cells.Item(212, 4).Value = 1424612653
$ws.Cells.Item(212, 5).Value = 0.475358
$ws.Cells.Item(212, 6).Value = 50.839662400000002

$ws.Cells.Item(213, 1).Value = 'Bat Cave'
$ws.Cells.Item(213, 2).Value = '5c Town Hall St, Blackburn BB2 1AG'
$ws.Cells.Item(213, 3).Value = 'BB2 1AG'
$ws.Cells.Item(213, 4).Value = 1254846356
$ws.Cells.Item(213, 5).Value = -2.4852764999999999
$ws.Cells.Item(213, 6).Value = 53.749734500000002

$ws.Cells.Item(214, 1).Value = 'Geek Retreat Blackburn'
$ws.Cells.Item(214, 2).Value = '50 - 54 Church St, Blackburn BB1 5AL'
$ws.Cells.Item(214, 3).Value = 'BB1 5AL'
$ws.Cells.Item(214, 4).Value = 7943622350
$ws.Cells.Item(214, 5).Value = -2.4814006000000002
$ws.Cells.Item(214, 6).Value = 53.747911999999999

$ws.Cells.Item(215, 1).Value = 'Tabletop Dominion'
$ws.Cells.Item(215, 2).Value = '14 Broadway, Accrington BB5 1EY'
$ws.Cells.Item(215, 3).Value = 'BB5 1EY'
$ws.Cells.Item(215, 4).Value = 7990388340
$ws.Cells.Item(215, 5).Value = -2.36503
$ws.Cells.Item(215, 6).Value = 53.753480000000003

# Match the author's final selection/view state.
$ws.Range("A204:F215").Select()
